# "solved the formcode.php and thus the comment section"
# - Restyle the header row (A1:F1) with a new, smaller, non-bold Arial font,
#   centered horizontally.
# - Widen column E to fit the new "Comment" values.
# - Add a new time-sheet entry in row 6 (start/end time, computed hours,
#   date, topic "Comment", and the new comment text).
# - Update the active selection to match the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row restyle (A1:F1): Arial 12, not bold, black, centered -------
$header = $ws.Range("A1:F1")
$header.Font.Name = "Arial"
$header.Font.Size = 12
$header.Font.Color = 0
$header.HorizontalAlignment = -4108   # xlCenter

# --- Column widths ----------------------------------------------------------
# (ColumnWidth is stored/rounded in sixths of a character width by this
# engine, so 14.75 is the closest input that round-trips to ~15.625/15.667.)
$ws.Columns.Item(5).ColumnWidth = 14.75

# --- New row 6: time entry ---------------------------------------------------
# Start time 13:00 / End time 13:54, copying formatting from the row above.
$ws.Cells.Item(5,1).Copy($ws.Cells.Item(6,1))
$ws.Cells.Item(6,1).Value = 0.54166666666666663

$ws.Cells.Item(5,2).Copy($ws.Cells.Item(6,2))
$ws.Cells.Item(6,2).Value = 0.57916666666666672

$ws.Cells.Item(5,3).Copy($ws.Cells.Item(6,3))
$ws.Cells.Item(6,3).Formula = "=B6-A6"

$ws.Cells.Item(5,4).Copy($ws.Cells.Item(6,4))
$ws.Cells.Item(6,4).Value = 41787

$ws.Cells.Item(5,5).Copy($ws.Cells.Item(6,5))
$ws.Cells.Item(6,5).Value = "Comment"

$ws.Cells.Item(5,6).Copy($ws.Cells.Item(6,6))
$ws.Cells.Item(6,6).Value = "Solved the php include stuff"

# --- Selection ---------------------------------------------------------------
$ws.Range("F9:F10").Select()

Write-Output "done"
